# Updated cryptos list - applies latest price/volume snapshot to the
# "cryptos" sheet, plus a ranking swap between Polygon and Chainlink
# (rows 15 and 16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing it to be stored as text.
# The Price column ("D") holds values that look numeric (e.g. "7.12",
# "243.05"); Excel's COM automation auto-converts such strings to real
# numbers when assigned straight to .Value on a General-formatted cell.
# Briefly switching the cell to the Text format before the assignment -
# then reverting the cell style back to Normal - keeps the stored type as
# text (matching the original inlineStr cells) without leaving any visible
# style/format change behind.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Cells.Item(2, 4) "42.042.21"
$ws.Range("E2").Value = "  -0.43%  "

# Row 3 - Ethereum
Set-TextValue $ws.Cells.Item(3, 4) "2.220.30"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5 - BNB
Set-TextValue $ws.Cells.Item(5, 4) "243.05"
$ws.Range("E5").Value = "  -1.70%  "

# Row 6 - XRP
Set-TextValue $ws.Cells.Item(6, 4) "0.627"
$ws.Range("E6").Value = "  -0.27%  "

# Row 7 - Solana
Set-TextValue $ws.Cells.Item(7, 4) "73.79"
$ws.Range("E7").Value = "  -1.01%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.15%  "

# Row 9 - Cardano
Set-TextValue $ws.Cells.Item(9, 4) "0.615"
$ws.Range("E9").Value = "  -0.60%  "

# Row 10 - Avalanche
Set-TextValue $ws.Cells.Item(10, 4) "43.90"
$ws.Range("E10").Value = "  +6.30%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.85%  "

# Row 12 - Polkadot
Set-TextValue $ws.Cells.Item(12, 4) "7.12"
$ws.Range("E12").Value = "  +0.28%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.69%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Cells.Item(14, 4) "2.546.56"

# Row 15 / 16 - Polygon and Chainlink swap ranking positions
Set-TextValue $ws.Cells.Item(15, 2) "Chainlink"
Set-TextValue $ws.Cells.Item(15, 3) "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Cells.Item(15, 4) "14.24"
$ws.Range("E15").Value = "  -1.87%  "

Set-TextValue $ws.Cells.Item(16, 2) "Polygon"
Set-TextValue $ws.Cells.Item(16, 3) "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue $ws.Cells.Item(16, 4) "0.844"
$ws.Range("E16").Value = "  -1.27%  "

# Row 17 - WrappedEther
Set-TextValue $ws.Cells.Item(17, 4) "2.212.32"
$ws.Range("E17").Value = "  -1.07%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Cells.Item(18, 4) "41.896.88"
$ws.Range("E18").Value = "  -0.34%  "

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  +12.66%  "

# Row 20 - Uniswap
Set-TextValue $ws.Cells.Item(20, 4) "6.18"
$ws.Range("E20").Value = "  +0.89%  "

# Row 21 - Litecoin
Set-TextValue $ws.Cells.Item(21, 4) "72.38"
$ws.Range("E21").Value = "  +0.70%  "

# Row 22 - InternetComputer(DFINITY)
$ws.Range("E22").Value = "  +34.36%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Cells.Item(23, 4) "229.58"
$ws.Range("E23").Value = "  -0.95%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  -7.97%  "

# Row 25 - Cosmos
Set-TextValue $ws.Cells.Item(25, 4) "11.54"
$ws.Range("E25").Value = "  +3.32%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.08%  "

# Row 27 - WEMIXToken
$ws.Range("E27").Value = "  +1.42%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -1.38%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +5.73%  "

# Row 30 - Monero
Set-TextValue $ws.Cells.Item(30, 4) "166.54"
$ws.Range("E30").Value = "  -1.83%  "

# Row 31 - EthereumClassic
Set-TextValue $ws.Cells.Item(31, 4) "20.62"
$ws.Range("E31").Value = "  -0.20%  "

# Row 32 - Filecoin
Set-TextValue $ws.Cells.Item(32, 4) "5.65"
$ws.Range("E32").Value = "  +15.38%  "

# Row 33 - Hedera
Set-TextValue $ws.Cells.Item(33, 4) "0.0799"
$ws.Range("E33").Value = "  -3.09%  "

# Row 34 - Stellar
Set-TextValue $ws.Cells.Item(34, 4) "0.125"
$ws.Range("E34").Value = "  -0.14%  "

# Row 35 - InjectiveProtocol
Set-TextValue $ws.Cells.Item(35, 4) "29.40"
$ws.Range("E35").Value = "  -2.65%  "

# Row 36 - Kaspa
$ws.Range("E36").Value = "  -3.95%  "

# Row 37 - RenderToken
Set-TextValue $ws.Cells.Item(37, 4) "4.30"
$ws.Range("E37").Value = "  -4.57%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.51%  "

# Row 39 - Celestia
Set-TextValue $ws.Cells.Item(39, 4) "13.00"
$ws.Range("E39").Value = "  -4.42%  "

# Row 40 - LidoDAOToken
Set-TextValue $ws.Cells.Item(40, 4) "2.14"
$ws.Range("E40").Value = "  -2.17%  "

# Row 41 - MultiversX
Set-TextValue $ws.Cells.Item(41, 4) "65.16"
$ws.Range("E41").Value = "  +4.70%  "

# Row 42 - THORChain
$ws.Range("E42").Value = "  -2.14%  "

# Row 43 - Algorand
$ws.Range("E43").Value = "  -1.39%  "

# Row 44 - FraxShare
Set-TextValue $ws.Cells.Item(44, 4) "8.71"
$ws.Range("E44").Value = "  +0.45%  "

# Row 45 - Aave
Set-TextValue $ws.Cells.Item(45, 4) "104.66"
$ws.Range("E45").Value = "  -3.58%  "

# Row 46 - Cronos
$ws.Range("E46").Value = "  +0.68%  "

# Row 47 - NEARProtocol
Set-TextValue $ws.Cells.Item(47, 4) "2.40"
$ws.Range("E47").Value = "  +5.00%  "

# Row 48 - ARBITRUM
$ws.Range("E48").Value = "  -0.64%  "

# Row 49 - TrustWalletToken
$ws.Range("E49").Value = "  -0.24%  "

# Row 50 - HuobiToken
$ws.Range("E50").Value = "  +0.82%  "

# Row 51 - RocketPoolETH
Set-TextValue $ws.Cells.Item(51, 4) "2.425.13"
$ws.Range("E51").Value = "  -1.31%  "
